$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 1, 46.85851166666667, 140.575535, 0.1419057303676978, 0.1419057303676978, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 2479.298081651376, 22313.68273486238, 0.05905286953281434, 0.05905286953281434),
    @(3, 1, 46.85851166666667, 140.575535, 0.1419057303676978, 0.1419057303676978, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 2214.535994779847, 19930.82395301862, 0.05274666493040277, 0.05274666493040278),
    @(3, 1, 46.85851166666667, 140.575535, 0.1419057303676978, 0.1419057303676978, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 1263.989952432752, 11375.90957189477, 0.0301061959044807, 0.0301061959044807),
    @(3, 1, 283.1772156666667, 849.531647, 0.8575703363889615, 0.8575703363889616, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 14982.99247240449, 134846.9322516404, 0.3568706426355616, 0.3568706426355617),
    @(3, 1, 283.1772156666667, 849.531647, 0.8575703363889615, 0.8575703363889616, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 13382.97173107758, 120446.7455796982, 0.3187607369382034, 0.3187607369382034),
    @(3, 1, 283.1772156666667, 849.531647, 0.8575703363889615, 0.8575703363889616, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 7638.594198354982, 68747.34778519484, 0.1819389568151965, 0.1819389568151965),
    @(3, 1, 0.1730073333333333, 0.519022, 0.000523933243340694, 0.000523933243340694, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 9.153870543226887, 82.384834889042, 0.0002180303880804036, 0.0002180303880804036),
    @(3, 1, 0.1730073333333333, 0.519022, 0.000523933243340694, 0.000523933243340694, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 8.176336665427776, 73.58702998884999, 0.0001947471124723623, 0.0001947471124723623),
    @(3, 1, 0.1730073333333333, 0.519022, 0.000523933243340694, 0.000523933243340694, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 4.666804882453778, 42.001243942084, 0.0001111557427879281, 0.0001111557427879281)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = 5 + $j
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

Write-Output "Done updating values"